# fix bug hard code in filter_date function
# Update the "last_edited_time" values (column D) on the LUY_KE_NGAY_LONG_XUYEN
# sheet so they reflect the corrected (non hard-coded) filter date/time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D29").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D30:D50").Value = "2024-07-04T09:48:00.000Z"
$ws.Range("D51:D65").Value = "2024-07-04T09:47:00.000Z"
$ws.Range("D66:D72").Value = "2024-07-04T09:48:00.000Z"
